$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HHC")

# Insert two new columns before column D (pushes existing D:K to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Apply number formats to the newly inserted columns to match the rest of the table
$ws.Range("D7:E7").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D38:E38").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D80:E80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D8:E35").NumberFormat = "#,##0"
$ws.Range("D39:E77").NumberFormat = "#,##0"
$ws.Range("D81:E102").NumberFormat = "#,##0"

# Set data values (new columns D,E; and corrected historical values for revised rows)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 464700
$ws.Range("E8").Value = 257200
$ws.Range("D9").Value = 307500
$ws.Range("E9").Value = 141400
$ws.Range("D10").Value = 157200
$ws.Range("E10").Value = 115800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 38200
$ws.Range("E15").Value = 31100
$ws.Range("D17").Value = 394400
$ws.Range("E17").Value = 214800
$ws.Range("D18").Value = 70300
$ws.Range("E18").Value = 42400
$ws.Range("D20").Value = 2400
$ws.Range("E20").Value = 10700
$ws.Range("D21").Value = 110800
$ws.Range("E21").Value = 84600
$ws.Range("D22").Value = 24800
$ws.Range("E22").Value = 21700
$ws.Range("D23").Value = 47800
$ws.Range("E23").Value = 31300
$ws.Range("D24").Value = 9900
$ws.Range("E24").Value = 7500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 37900
$ws.Range("E26").Value = 23800
$ws.Range("D27").Value = 37300
$ws.Range("E27").Value = 23400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2400
$ws.Range("E32").Value = -10700
$ws.Range("D33").Value = 37300
$ws.Range("E33").Value = 23400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 37300
$ws.Range("E35").Value = 23400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 499700
$ws.Range("E41").Value = 454100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 266900
$ws.Range("E43").Value = 293700
$ws.Range("F43").Value = 279500
$ws.Range("G43").Value = 254600
$ws.Range("H43").Value = 393500
$ws.Range("I43").Value = 489200
$ws.Range("J43").Value = 589200
$ws.Range("K43").Value = 565700
$ws.Range("L43").Value = 488100
$ws.Range("M43").Value = 466200
$ws.Range("D44").Value = 198400
$ws.Range("E44").Value = 78700
$ws.Range("D45").Value = 17000
$ws.Range("E45").Value = 16100
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 107000
$ws.Range("E47").Value = 147900
$ws.Range("D48").Value = 5797900
$ws.Range("E48").Value = 5917300
$ws.Range("D49").Value = 59800
$ws.Range("E49").Value = 61100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 320300
$ws.Range("E52").Value = 254300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 7355800
$ws.Range("E54").Value = 7296900
$ws.Range("D57").Value = 38700
$ws.Range("E57").Value = 32600
$ws.Range("F57").Value = 21900
$ws.Range("G57").Value = 26300
$ws.Range("H57").Value = 35900
$ws.Range("I57").Value = 33500
$ws.Range("J57").Value = 33800
$ws.Range("K57").Value = 233800
$ws.Range("L57").Value = 241000
$ws.Range("M57").Value = 230700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 635100
$ws.Range("E59").Value = 584400
$ws.Range("F59").Value = 585000
$ws.Range("G59").Value = 497200
$ws.Range("H59").Value = 376600
$ws.Range("I59").Value = 324000
$ws.Range("J59").Value = 334100
$ws.Range("K59").Value = 174800
$ws.Range("L59").Value = 216100
$ws.Range("M59").Value = 233200
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 3181200
$ws.Range("E61").Value = 3261200
$ws.Range("D62").Value = 157200
$ws.Range("E62").Value = 148800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 4223600
$ws.Range("E66").Value = 4191600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -120300
$ws.Range("E72").Value = -157600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3132200
$ws.Range("E76").Value = 3105300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 37300
$ws.Range("E81").Value = 23400
$ws.Range("D83").Value = 38200
$ws.Range("E83").Value = 31600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 350300
$ws.Range("E89").Value = -28000
$ws.Range("F89").Value = -25400
$ws.Range("G89").Value = -86400
$ws.Range("H89").Value = 343400
$ws.Range("I89").Value = 73300
$ws.Range("J89").Value = -39900
$ws.Range("K89").Value = -57800
$ws.Range("L89").Value = 277200
$ws.Range("M89").Value = -81500
$ws.Range("D91").Value = 200
$ws.Range("E91").Value = -1700
$ws.Range("F91").Value = -1700
$ws.Range("G91").Value = -1300
$ws.Range("H91").Value = -1000
$ws.Range("I91").Value = -1100
$ws.Range("J91").Value = -2300
$ws.Range("K91").Value = -237900
$ws.Range("L91").Value = -109900
$ws.Range("M91").Value = -95300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -187600
$ws.Range("E94").Value = -253000
$ws.Range("F94").Value = -301200
$ws.Range("G94").Value = -99900
$ws.Range("H94").Value = -70600
$ws.Range("I94").Value = -70000
$ws.Range("J94").Value = -90000
$ws.Range("K94").Value = -85000
$ws.Range("L94").Value = -104500
$ws.Range("M94").Value = -95900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -51000
$ws.Range("E100").Value = 157200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 111700
$ws.Range("E102").Value = -123800
$ws.Range("F102").Value = -28600
$ws.Range("G102").Value = -199400
$ws.Range("H102").Value = 148500
$ws.Range("I102").Value = 52500
$ws.Range("J102").Value = 118600
$ws.Range("K102").Value = -124000
$ws.Range("L102").Value = 12500
$ws.Range("M102").Value = -17800
